{"js": "// 1) Move the \"_GoBack\" bookmark to the very start of the document\n//    (it currently sits inside the last paragraph, splitting a run in two).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst body = context.document.body;\nconst firstParagraph = body.paragraphs.getFirst();\nconst startOfDoc = firstParagraph.getRange(\"Start\");\nstartOfDoc.insertBookmark(\"_GoBack\");\n\n// 2) Un-split the \"- As hip\u00f3teses a serem testadas na ANOVA;\" sentence: it is\n//    currently stored as a run-per-syllable; replacing the matched range with\n//    the same text collapses it back into a single run.\nconst hipoteses = body.search(\"- As hip\u00f3teses a serem testadas na ANOVA;\", { matchCase: false });\nhipoteses.load(\"items\");\nawait context.sync();\nhipoteses.items[0].insertText(\"- As hip\u00f3teses a serem testadas na ANOVA;\", \"Replace\");\n\n// 3) Likewise merge the trailing \". Apresente a tabela do teste com as\n//    m\u00e9dias.\" sentence (previously split around the bookmark) into one run.\nconst medias = body.search(\". Apresente a tabela do teste com as m\u00e9dias.\", { matchCase: false });\nmedias.load(\"items\");\nawait context.sync();\nmedias.items[0].insertText(\". Apresente a tabela do teste com as m\u00e9dias.\", \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Move the \"_GoBack\" bookmark to the very start of the document (it\n#    currently sits inside the last paragraph, splitting a run in two).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n$startOfDoc = $d.Range(0, 0)\n$d.Bookmarks.Add(\"_GoBack\", $startOfDoc)\n\n# 2) Un-split the \"- As hip\u00f3teses a serem testadas na ANOVA;\" sentence: it is\n#    currently stored as a run-per-syllable; a Find/Replace over the exact\n#    matched text collapses it back into a single run.\n$rngHipoteses = $d.Content\n$rngHipoteses.Find.Execute(\"- As hip\u00f3teses a serem testadas na ANOVA;\", $false, $false, $false, $false, $false, $true, 1, $false, \"- As hip\u00f3teses a serem testadas na ANOVA;\", 2)\n\n# 3) Likewise merge the trailing \". Apresente a tabela do teste com as\n#    m\u00e9dias.\" sentence (previously split around the bookmark) into one run.\n$rngMedias = $d.Content\n$rngMedias.Find.Execute(\". Apresente a tabela do teste com as m\u00e9dias.\", $false, $false, $false, $false, $false, $true, 1, $false, \". Apresente a tabela do teste com as m\u00e9dias.\", 2)\n"}
